# Apply "Add data for 2022-12-21" changes: the partial-month column (B) now
# covers "through December 13" instead of "through December 12", and a
# number of neighborhood/month cells pick up incremented or newly-added
# carjacking counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Sheet title / header text -------------------------------------------
$ws.Name = "Through 2022-12-13"
$ws.Range("B1").Value = "December 2022 (through December 13)"

# --- Row 2  (Garfield Park) ------------------------------------------------
$ws.Range("Z2").Value = 9
$ws.Range("BJ2").Value = 5
$ws.Range("BV2").Value = 2

# --- Row 4  (Chatham) ------------------------------------------------------
$ws.Range("BV4").Value = 2

# --- Row 8  (Belmont Cragin) ------------------------------------------------
$ws.Range("Z8").Value = 1

# --- Row 9  (South Shore) ---------------------------------------------------
$ws.Range("Z9").Value = 3
$ws.Range("BV9").Value = 3

# --- Row 10 (Grand Boulevard) ------------------------------------------------
$ws.Range("N10").Value = 4

# --- Row 12 (New City) -------------------------------------------------------
$ws.Range("AL12").Value = 2

# --- Row 14 (Austin) ----------------------------------------------------------
$ws.Range("BJ14").Value = 3

# --- Row 18 (Chicago Lawn) -----------------------------------------------------
$ws.Range("AX18").Value = 1

# --- Row 20 (North Lawndale) ----------------------------------------------------
$ws.Range("N20").Value = 4
$ws.Range("Z20").Value = 9

# --- Row 21 (West Town) ----------------------------------------------------------
$ws.Range("B21").Value = 3
$ws.Range("CH21").Value = 1

# --- Row 28 (Auburn Gresham) -------------------------------------------------------
$ws.Range("B28").Value = 2

# --- Row 35 (Wicker Park) ------------------------------------------------------------
$ws.Range("AL35").Value = 1

# --- Row 39 (Avalon Park) --------------------------------------------------------------
$ws.Range("AL39").Value = 1

# --- Row 51 (Archer Heights) -------------------------------------------------------------
$ws.Range("N51").Value = 1

# --- Row 55 (Bucktown) ---------------------------------------------------------------------
$ws.Range("B55").Value = 2

# --- Row 57 (Chinatown) ---------------------------------------------------------------------
$ws.Range("N57").Value = 4

# --- Row 59 (Douglas) -----------------------------------------------------------------------
$ws.Range("AL59").Value = 1

# --- Row 75 (Lincoln Park) -------------------------------------------------------------------
$ws.Range("N75").Value = 2
$ws.Range("Z75").Value = 1

# --- Row 84 (North Center) -------------------------------------------------------------------
$ws.Range("BJ84").Value = 2

# --- Row 92 (Rogers Park) --------------------------------------------------------------------
$ws.Range("Z92").Value = 1
